$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Semana del 11 de marzo (columna E, fila de "Camila"): actualizar el total de citas creadas
$ws.Range("E3").Value = 1528

# Dejar la celda E4 como seleccion activa, como quedo tras editar la hoja
$ws.Range("E4").Select()
